$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "74.784.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.838.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.07%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "599.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.36%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.77%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -5.25%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.837.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.63%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.25%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.369"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.02%  "

# Row 13 - Toncoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.71%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.372.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.68%  "

# Row 15 - WrappedBTC
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.851.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16 - row16(was ShibaInu -> Avalanche)
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.89%  "

# Row 17 - row17(was Avalanche -> ShibaInu)
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000187"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "

# Row 18 - WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.845.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.41%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.78%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "

# Row 22 - SuiNetwork
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.59%  "

# Row 23 - Polkadot
$ws.Range("E23").Value = "  +2.14%  "

# Row 24 - LEO
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.76%  "

# Row 25 - Dai
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "

# Row 26 - Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.56%  "

# Row 27 - row27(was NEARProtocol -> WrappedeETH)
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.983.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.82%  "

# Row 28 - row28(was WrappedeETH -> NEARProtocol)
$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "

# Row 29 - Aptos
$ws.Range("E29").Value = "  +4.49%  "

# Row 30 - PEPE
$ws.Range("E30").Value = "  +10.06%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +0.17%  "

# Row 32 - Bittensor
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "530.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.61%  "

# Row 33 - Fetch.AI
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.55%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "

# Row 35 - PancakeSwap
$ws.Range("E35").Value = "  +5.27%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  -0.05%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +2.25%  "

# Row 38 - EthereumClassic
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.56%  "

# Row 39 - Monero
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.40%  "

# Row 40 - WhiteBITCoin
$ws.Range("E40").Value = "  -0.46%  "

# Row 41 - Aave
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "184.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +22.33%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  +3.50%  "

# Row 44 - PolygonEcosystemToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.341"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.47%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +1.56%  "

# Row 46 - ImmutableX
$ws.Range("E46").Value = "  +6.23%  "

# Row 47 - OKB
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.76%  "

# Row 48 - dogwifhat
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.26%  "

# Row 49 - Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0863"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.71%  "

# Row 50 - ARBITRUM
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.572"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.96%  "

# Row 51 - Filecoin
$ws.Range("E51").Value = "  +4.01%  "
